# Apply the "[base commands] - [assertMatch(text,regex)]: *NEW* command" edit.
#
# The workbook keeps a hidden "#system" sheet that is a flat per-category
# lookup table: each column holds the header (category name) in row 1 and the
# command list for that category below it. Adding the new assertMatch/openFile
# commands means inserting rows into two of those column-lists (shifting the
# remaining entries down), and the loss of "tn.5250" as a selectable "target"
# means both a row removal from the target list and (because the tn.5250
# command sub-list column is no longer used) an entire column removal further
# right in the sheet.
#
# We use plain cell-value shifting (rather than Range.Insert/Delete) for the
# row-level edits because Range.Insert/Delete on a single-column range here
# shifts the *whole row* across every column, not just the target column.
# EntireColumn.Delete() is safe/columnar, so that one is used directly for
# the column removal.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

function Insert-ColumnValue {
    param($ws, $colIndex, $insertRow, $lastRow, $newValue)
    # Shift existing values in column $colIndex down by one row, starting
    # from the bottom, to open up a slot at $insertRow. $lastRow is the last
    # currently-populated row in that column (before the insert).
    for ($r = $lastRow + 1; $r -gt $insertRow; $r--) {
        $srcVal = $ws.Cells.Item($r - 1, $colIndex).Value()
        $ws.Cells.Item($r, $colIndex).Value = $srcVal
    }
    $ws.Cells.Item($insertRow, $colIndex).Value = $newValue
}

function Delete-ColumnValue {
    param($ws, $colIndex, $deleteRow, $lastRow)
    # Shift values above $deleteRow up by one row (closing the gap left by
    # removing the entry at $deleteRow), then blank out the now-unused last
    # row. $lastRow is the last currently-populated row in that column
    # (before the delete).
    for ($r = $deleteRow; $r -lt $lastRow; $r++) {
        $srcVal = $ws.Cells.Item($r + 1, $colIndex).Value()
        $ws.Cells.Item($r, $colIndex).Value = $srcVal
    }
    $ws.Cells.Item($lastRow, $colIndex).ClearContents()
}

# --- 1. "base" command list (column F): new assertMatch(text,regex) command,
#         inserted alphabetically between assertEqual and assertNotContain.
Insert-ColumnValue $ws 6 11 42 "assertMatch(text,regex)"

# --- 2. "external" command list (column J): new openFile(filePath) command,
#         inserted alphabetically at the top of the list.
Insert-ColumnValue $ws 10 2 6 "openFile(filePath)"

# --- 3. "target" category list (column A): tn.5250 is no longer offered as
#         a top-level target, so remove it from the list.
Delete-ColumnValue $ws 1 27 33

# --- 4. The tn.5250 command sub-list (old column AA) is now unused, so the
#         whole column is removed, shifting web/webalert/webcookie/ws/
#         ws.async/xml one column to the left (AB..AG -> AA..AF).
$ws.Range("AA1").EntireColumn.Delete()

# --- 5. Update the defined names whose ranges moved because of the edits
#         above (this engine does not auto-adjust Name.RefersTo on
#         insert/delete the way desktop Excel does). Note: the "tn.5250"
#         name intentionally keeps its old, now-stale reference -- it is no
#         longer reachable from the "target" list, so nothing repoints it.
$wb.Names.Item("base").RefersTo       = "='#system'!`$F`$2:`$F`$45"
$wb.Names.Item("external").RefersTo   = "='#system'!`$J`$2:`$J`$7"
$wb.Names.Item("target").RefersTo     = "='#system'!`$A`$2:`$A`$32"
$wb.Names.Item("web").RefersTo        = "='#system'!`$AA`$2:`$AA`$151"
$wb.Names.Item("webalert").RefersTo   = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("webcookie").RefersTo  = "='#system'!`$AC`$2:`$AC`$10"
$wb.Names.Item("ws").RefersTo         = "='#system'!`$AD`$2:`$AD`$17"
$wb.Names.Item("ws.async").RefersTo   = "='#system'!`$AE`$2:`$AE`$8"
$wb.Names.Item("xml").RefersTo        = "='#system'!`$AF`$2:`$AF`$27"
